$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1514.7858
$ws.Range("I32").Value = 1566.6666
$ws.Range("J32").Value = 1500.6364
$ws.Range("K32").Value = 1566.6666
$ws.Range("L32").Value = 1500.6364
$ws.Range("M32").Value = -1240.6666
$ws.Range("N32").Value = -2152.6364

$ws.Range("H69").Value = 3837.375
$ws.Range("I69").Value = 3600
$ws.Range("J69").Value = 4074.75
$ws.Range("K69").Value = 10800
$ws.Range("L69").Value = 12224.25
$ws.Range("M69").Value = -9926
$ws.Range("N69").Value = -13972.25

$ws.Range("H70").Value = 1247.8572
$ws.Range("I70").Value = 1083.25
$ws.Range("J70").Value = 1467.3334
$ws.Range("K70").Value = 3249.75
$ws.Range("L70").Value = 4402.0002
$ws.Range("M70").Value = -2979.75
$ws.Range("N70").Value = -4942.0002

$ws.Range("H72").Value = 3837.375
$ws.Range("I72").Value = 3600
$ws.Range("J72").Value = 4074.75
$ws.Range("K72").Value = 32400
$ws.Range("L72").Value = 36672.75
$ws.Range("M72").Value = -28032
$ws.Range("N72").Value = -45408.75

$ws.Range("H73").Value = 1247.8572
$ws.Range("I73").Value = 1083.25
$ws.Range("J73").Value = 1467.3334
$ws.Range("K73").Value = 3249.75
$ws.Range("L73").Value = 4402.0002
$ws.Range("M73").Value = -2313.75
$ws.Range("N73").Value = -6274.0002

$ws.Range("H80").Value = 20202462
$ws.Range("I80").Value = 30303456
$ws.Range("J80").Value = 475.36365
$ws.Range("K80").Value = 90910368
$ws.Range("L80").Value = 1426.09095
$ws.Range("M80").Value = -90909370
$ws.Range("N80").Value = -3422.09095

$ws.Range("H83").Value = 20202462
$ws.Range("I83").Value = 30303456
$ws.Range("J83").Value = 475.36365
$ws.Range("K83").Value = 272731104
$ws.Range("L83").Value = 4278.27285
$ws.Range("M83").Value = -272726112
$ws.Range("N83").Value = -14262.27285

$ws.Range("H123").Value = 39996.668
$ws.Range("J123").Value = 39996.668
$ws.Range("L123").Value = 39996.668
$ws.Range("N123").Value = -49796.668

$ws.Range("H124").Value = 51623.332
$ws.Range("J124").Value = 51623.332
$ws.Range("L124").Value = 51623.332
$ws.Range("N124").Value = -61443.332

$ws.Range("H126").Value = 47752
$ws.Range("J126").Value = 47752
$ws.Range("L126").Value = 47752
$ws.Range("N126").Value = -57632

$ws.Range("H128").Value = 52251.75
$ws.Range("J128").Value = 52251.75
$ws.Range("L128").Value = 52251.75
$ws.Range("N128").Value = -62211.75

$ws.Range("H130").Value = 56986.668
$ws.Range("J130").Value = 56986.668
$ws.Range("L130").Value = 56986.668
$ws.Range("N130").Value = -67026.66800000001

$ws.Range("H141").Value = 3385.1924
$ws.Range("I141").Value = 1985.7142
$ws.Range("J141").Value = 9263
$ws.Range("K141").Value = 5957.142599999999
$ws.Range("L141").Value = 27789
$ws.Range("M141").Value = -777.1425999999992
$ws.Range("N141").Value = -38149

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2822.2144
$ws.Range("I74").Value = 2204
$ws.Range("J74").Value = 2990.818
$ws.Range("K74").Value = 2204
$ws.Range("L74").Value = 2990.818
$ws.Range("M74").Value = -1330
$ws.Range("N74").Value = -4738.818

$ws.Range("H77").Value = 2822.2144
$ws.Range("I77").Value = 2204
$ws.Range("J77").Value = 2990.818
$ws.Range("K77").Value = 11020
$ws.Range("L77").Value = 14954.09
$ws.Range("M77").Value = -6652
$ws.Range("N77").Value = -23690.09

$ws.Range("H130").Value = 43750
$ws.Range("J130").Value = 43750
$ws.Range("L130").Value = 43750
$ws.Range("N130").Value = -53790

$ws.Range("H131").Value = 49707
$ws.Range("J131").Value = 49707
$ws.Range("L131").Value = 49707
$ws.Range("N131").Value = -59787

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 49233.332
$ws.Range("J130").Value = 49233.332
$ws.Range("L130").Value = 49233.332
$ws.Range("N130").Value = -59273.332

$ws.Range("H134").Value = 4853.636
$ws.Range("I134").Value = 4632.3335
$ws.Range("J134").Value = 5119.2
$ws.Range("K134").Value = 13897.0005
$ws.Range("L134").Value = 15357.6
$ws.Range("M134").Value = -11362.0005
$ws.Range("N134").Value = -20427.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 204549.16
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 215262.27
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 215262.27
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -215486.27

$ws.Range("H31").Value = 8572.879999999999
$ws.Range("I31").Value = 4501.4287
$ws.Range("K31").Value = 4501.4287
$ws.Range("M31").Value = -4206.4287

$ws.Range("H34").Value = 8572.879999999999
$ws.Range("I34").Value = 4501.4287
$ws.Range("K34").Value = 4501.4287
$ws.Range("M34").Value = -4299.4287

$ws.Range("H133").Value = 24893.666
$ws.Range("J133").Value = 24893.666
$ws.Range("L133").Value = 24893.666
$ws.Range("N133").Value = -29953.666

$ws.Range("H134").Value = 4139.5557
$ws.Range("I134").Value = 1300
$ws.Range("J134").Value = 5946.5454
$ws.Range("K134").Value = 3900
$ws.Range("L134").Value = 17839.6362
$ws.Range("M134").Value = -1365
$ws.Range("N134").Value = -22909.6362

$ws.Range("H137").Value = 25388.21
$ws.Range("J137").Value = 25388.21
$ws.Range("L137").Value = 25388.21
$ws.Range("N137").Value = -35588.21

$ws.Range("H141").Value = 19163.666
$ws.Range("J141").Value = 19495.5
$ws.Range("L141").Value = 19495.5
$ws.Range("N141").Value = -29855.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 696.2143
$ws.Range("I5").Value = 439.76
$ws.Range("J5").Value = 2833.3333
$ws.Range("K5").Value = 1319.28
$ws.Range("L5").Value = 8499.999899999999
$ws.Range("M5").Value = -1207.28
$ws.Range("N5").Value = -8723.999899999999

$ws.Range("H57").Value = 12884.308
$ws.Range("I57").Value = 24099.4
$ws.Range("J57").Value = 5874.875
$ws.Range("K57").Value = 72298.20000000001
$ws.Range("L57").Value = 17624.625
$ws.Range("M57").Value = -71739.20000000001
$ws.Range("N57").Value = -18742.625

$ws.Range("H58").Value = 2042896
$ws.Range("J58").Value = 2552620.2
$ws.Range("L58").Value = 7657860.600000001
$ws.Range("N58").Value = -7658116.600000001

$ws.Range("H135").Value = 696.2143
$ws.Range("I135").Value = 439.76
$ws.Range("J135").Value = 2833.3333
$ws.Range("K135").Value = 3957.84
$ws.Range("L135").Value = 25499.9997
$ws.Range("M135").Value = -1422.84
$ws.Range("N135").Value = -30569.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 13884.615
$ws.Range("J5").Value = 13884.615
$ws.Range("L5").Value = 13884.615
$ws.Range("N5").Value = -14108.615

$ws.Range("H80").Value = 211997.92
$ws.Range("I80").Value = 388800
$ws.Range("J80").Value = 3050
$ws.Range("K80").Value = 388800
$ws.Range("L80").Value = 3050
$ws.Range("M80").Value = -387802
$ws.Range("N80").Value = -5046

$ws.Range("H83").Value = 211997.92
$ws.Range("I83").Value = 388800
$ws.Range("J83").Value = 3050
$ws.Range("K83").Value = 1944000
$ws.Range("L83").Value = 15250
$ws.Range("M83").Value = -1939008
$ws.Range("N83").Value = -25234

$ws.Range("H130").Value = 51398.4
$ws.Range("J130").Value = 51398.4
$ws.Range("L130").Value = 51398.4
$ws.Range("N130").Value = -61438.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 7916.6113
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 24599.8
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 24599.8
$ws.Range("M2").Value = -1388
$ws.Range("N2").Value = -24823.8

$ws.Range("H119").Value = 45830.668
$ws.Range("J119").Value = 45830.668
$ws.Range("L119").Value = 45830.668
$ws.Range("N119").Value = -55506.668

$ws.Range("H132").Value = 3368.3635
$ws.Range("I132").Value = 2918.5144
$ws.Range("J132").Value = 4155.6
$ws.Range("K132").Value = 8755.5432
$ws.Range("L132").Value = 12466.8
$ws.Range("M132").Value = -6225.5432
$ws.Range("N132").Value = -17526.8

$ws.Range("H137").Value = 59138.168
$ws.Range("J137").Value = 59138.168
$ws.Range("L137").Value = 59138.168
$ws.Range("N137").Value = -69338.16800000001

$ws.Range("H139").Value = 57731.75
$ws.Range("J139").Value = 57731.75
$ws.Range("L139").Value = 57731.75
$ws.Range("N139").Value = -68011.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1542.8286
$ws.Range("I132").Value = 1230.1666
$ws.Range("J132").Value = 3418.8
$ws.Range("K132").Value = 3690.4998
$ws.Range("L132").Value = 10256.4
$ws.Range("M132").Value = -1160.4998
$ws.Range("N132").Value = -15316.4

$ws.Range("H136").Value = 18989.654
$ws.Range("I136").Value = 32387.5
$ws.Range("K136").Value = 97162.5
$ws.Range("M136").Value = -94612.5

$ws.Range("H138").Value = 39666.668
$ws.Range("J138").Value = 39666.668
$ws.Range("L138").Value = 39666.668
$ws.Range("N138").Value = -49946.668
